$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45108
$ws.Range("B2").Value = 906
$ws.Range("C2").Value = 6

$ws.Range("C3").Select()
